$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Time Period" range text (row 6, columns B:C)
$ws.Range("B6:C6").Value = "2018-03-02:2024-03-15"

# Update the "Update" date text (row 8, columns B:C) - keep as text, not auto-converted to a date
$ws.Range("B8:C8").Value = "'2024-03-15"

# Append new weekly data rows (312-320), copying formatting from the last existing row (311)
$ws.Range("A311:C311").Copy()
$ws.Range("A312:C320").PasteSpecial(-4122)

$newData = @(
    @(312, 45310, 76.25, 83),
    @(313, 45317, 76.84, 83.52),
    @(314, 45324, 76.5, 83.59),
    @(315, 45331, 76.78, 83.82),
    @(316, 45338, 76.349999999999994, 83.99),
    @(317, 45345, 75.61, 83.57),
    @(318, 45352, 75.17, 83.32),
    @(319, 45359, 75.62, 83.09),
    @(320, 45366, 76.17, 82.56)
)

foreach ($row in $newData) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

$excel.CutCopyMode = 0
